$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "goal" ramp labels in row 6 (Week 3 weekday/weekend goal markers)
$ws.Range("Z6").Value = "4 hours"
$ws.Range("AC6").Value = "3 hours"
$ws.Range("AE6").Value = "2 hours"

# Update the "running goal" milestone labels in row 8
$ws.Range("AC8").Value = "53 hours"
$ws.Range("AD8").Value = "58 hours"

# Add a new time-log entry on row 25 (date + start time)
$ws.Range("B25").Value = 45215
$ws.Range("C25").Value = 0.80555555555555547

# Move selection to reflect the new active cell / scrolled view
$ws.Range("D25").Select()
